$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume/Number and report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Simple numeric updates (style/type unchanged) ---
$numericUpdates = @{
    "C16" = 5
    "C17" = 9
    "C18" = 8
    "C20" = 2
    "C21" = 34
    "C24" = 34
    "C25" = 21
    "C26" = 18
    "D16" = 6
    "D17" = 5
    "D18" = 5
    "D19" = 14
    "D21" = 32
    "D24" = 47
    "D25" = 33
    "D26" = 12
    "E16" = -16.666666666666
    "E17" = 80
    "E18" = 60
    "E19" = -28.571428571428
    "E21" = 6.25
    "E24" = -27.659574468085
    "E25" = -36.363636363636
    "E26" = 50
    "E27" = -100
    "F16" = 13
    "F17" = 20
    "F18" = 18
    "F19" = 50
    "F20" = 11
    "F21" = 114
    "F22" = 5
    "F24" = 112
    "F25" = 62
    "F26" = 42
    "F28" = 6
    "G16" = 14
    "G17" = 23
    "G18" = 21
    "G19" = 48
    "G20" = 10
    "G21" = 116
    "G24" = 177
    "G26" = 55
    "G27" = 2
    "G28" = 4
    "G31" = 2
    "H16" = -7.142857142857
    "H17" = -13.043478260869
    "H18" = -14.285714285714
    "H19" = 4.166666666666
    "H20" = 10
    "H21" = -1.724137931034
    "H22" = 66.666666666666
    "H24" = -36.723163841807
    "H25" = -44.144144144144
    "H26" = -23.636363636363
    "H27" = 0
    "H28" = 50
    "I16" = 23
    "I17" = 37
    "I18" = 38
    "I19" = 87
    "I20" = 22
    "I21" = 211
    "I22" = 6
    "I24" = 269
    "I25" = 157
    "I26" = 82
    "I28" = 7
    "J16" = 38
    "J17" = 39
    "J18" = 41
    "J19" = 124
    "J20" = 20
    "J21" = 262
    "J24" = 361
    "J25" = 225
    "J26" = 92
    "J27" = 3
    "J31" = 4
    "K16" = -39.473684210526
    "K17" = -5.128205128205
    "K18" = -7.317073170731
    "K19" = -29.838709677419
    "K20" = 10
    "K21" = -19.465648854961
    "K22" = 0
    "K24" = -25.484764542936
    "K25" = -30.222222222222
    "K26" = -10.869565217391
    "K27" = 33.333333333333
    "K28" = -30
    "L15" = -75
    "L16" = -46.511627906976
    "L17" = 23.333333333333
    "L18" = 11.764705882352
    "L19" = -26.271186440678
    "L20" = -26.666666666666
    "L21" = -19.771863117870
    "L22" = -50
    "L24" = -0.738007380073
    "L25" = 28.688524590163
    "L26" = -11.827956989247
    "L27" = -60
    "L28" = -58.823529411764
    "M15" = 0
    "M16" = -37.837837837837
    "M17" = 117.647058823529
    "M18" = -5
    "M19" = 29.850746268656
    "M20" = -43.589743589743
    "M21" = 3.940886699507
    "M22" = -33.333333333333
    "M24" = 140.178571428571
    "M26" = -6.818181818181
    "N16" = -89.252336448598
    "N17" = -7.5
    "N18" = -84.615384615384
    "N19" = -37.410071942446
    "N20" = -93.491124260355
    "N21" = -78.491335372069
}
foreach ($ref in $numericUpdates.Keys) {
    $ws.Range($ref).Value = $numericUpdates[$ref]
}

# --- Cells changing from text ("0" / "***.*") to numeric values ---
# First copy formatting+type from a stable style-14/15 donor cell (row 14, untouched by this edit),
# then overwrite with the required numeric value.
$numericSwaps = @{
    "D20" = 2
    "D31" = 1
    "E20" = 0
    "E31" = -100
}
$numericDonor = @{
    "D20" = "F14"
    "E20" = "M14"
    "D31" = "I14"
    "E31" = "N14"
}
foreach ($ref in $numericSwaps.Keys) {
    $ws.Range($numericDonor[$ref]).Copy($ws.Range($ref))
    $ws.Range($ref).Value = $numericSwaps[$ref]
}

# --- Cells changing from numeric values to text ("0" / "***.*") ---
# Copy both value and formatting directly from a stable donor cell (row 14) that already
# holds the exact shared-string text we need.
$stringDonor = @{
    "C15" = "C14"
    "C27" = "C14"
    "D22" = "C14"
    "D33" = "C14"
    "E22" = "E14"
    "E33" = "E14"
}
foreach ($ref in $stringDonor.Keys) {
    $ws.Range($stringDonor[$ref]).Copy($ws.Range($ref))
}
